# update crisis year run
$wb = $excel.ActiveWorkbook

$wsLCOE = $wb.Worksheets.Item("LCOE")

# Updated crisis_year figures (row 6) on the LCOE sheet
$wsLCOE.Range("B6").Value = 241.6948431289498
$wsLCOE.Range("C6").Value = 1336.0353828516941
$wsLCOE.Range("E6").Value = 14550789.14001311

# Move the active tab/selection from LCOE_Sensitivities to LCOE
$wsLCOE.Activate() | Out-Null
$wsLCOE.Range("L12").Select() | Out-Null
